# Applies the Sat Oct 26 03:10:34 UTC 2024 GitHub Actions refresh of the
# cryptos list: updated prices/volumes and re-ranked several coins.
#
# All target cells hold plain text values (prices/percentages are stored as
# text, not numbers), so we force text interpretation via NumberFormat = "@"
# before assigning, then ClearFormats() to drop the temporary text format and
# keep the cell style identical to the original (no explicit style index).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
Set-TextValue "D2" '66.721.08'
Set-TextValue "E2" '  -2.03%  '
# Row 3
Set-TextValue "D3" '2.440.40'
Set-TextValue "E3" '  -3.36%  '
# Row 4
Set-TextValue "D4" '1.00'
Set-TextValue "E4" '  -0.02%  '
# Row 5
Set-TextValue "D5" '576.68'
Set-TextValue "E5" '  -3.38%  '
# Row 6
Set-TextValue "D6" '164.24'
Set-TextValue "E6" '  -6.58%  '
# Row 8
Set-TextValue "D8" '0.509'
Set-TextValue "E8" '  -3.88%  '
# Row 9
Set-TextValue "D9" '2.440.88'
Set-TextValue "E9" '  -3.30%  '
# Row 10
Set-TextValue "E10" '  -5.08%  '
# Row 11
Set-TextValue "D11" '0.163'
Set-TextValue "E11" '  -0.94%  '
# Row 12
Set-TextValue "B12" 'Toncoin'
Set-TextValue "C12" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D12" '4.85'
Set-TextValue "E12" '  -5.08%  '
# Row 13
Set-TextValue "B13" 'Cardano'
Set-TextValue "C13" 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue "D13" '0.330'
Set-TextValue "E13" '  -4.19%  '
# Row 14
Set-TextValue "B14" 'WrappedliquidstakedEther2.0'
Set-TextValue "C14" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue "D14" '2.888.19'
Set-TextValue "E14" '  -3.23%  '
# Row 15
Set-TextValue "B15" 'Avalanche'
Set-TextValue "C15" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D15" '25.13'
Set-TextValue "E15" '  -5.62%  '
# Row 16
Set-TextValue "D16" '66.698.63'
Set-TextValue "E16" '  -1.72%  '
# Row 17
Set-TextValue "E17" '  -6.67%  '
# Row 18
Set-TextValue "D18" '2.478.28'
Set-TextValue "E18" '  -2.34%  '
# Row 19
Set-TextValue "B19" 'Chainlink'
Set-TextValue "C19" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D19" '11.24'
Set-TextValue "E19" '  -6.87%  '
# Row 20
Set-TextValue "B20" 'Uniswap'
Set-TextValue "C20" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue "D20" '7.68'
Set-TextValue "E20" '  -5.32%  '
# Row 21
Set-TextValue "B21" 'BitcoinCash'
Set-TextValue "C21" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue "D21" '352.82'
Set-TextValue "E21" '  -3.54%  '
# Row 22
Set-TextValue "D22" '4.02'
Set-TextValue "E22" '  -3.67%  '
# Row 23
Set-TextValue "E23" '  +0.04%  '
# Row 24
Set-TextValue "D24" '69.05'
Set-TextValue "E24" '  -2.98%  '
# Row 25
Set-TextValue "D25" '4.21'
Set-TextValue "E25" '  -9.82%  '
# Row 27
Set-TextValue "D27" '8.86'
Set-TextValue "E27" '  -12.90%  '
# Row 28
Set-TextValue "D28" '0.999'
Set-TextValue "E28" '  -0.22%  '
# Row 30
Set-TextValue "D30" '0.0₃0893'
Set-TextValue "E30" '  -9.19%  '
# Row 31
Set-TextValue "D31" '503.06'
Set-TextValue "E31" '  -5.54%  '
# Row 32
Set-TextValue "D32" '7.78'
Set-TextValue "E32" '  -6.67%  '
# Row 33
Set-TextValue "D33" '1.76'
Set-TextValue "E33" '  -7.99%  '
# Row 34
Set-TextValue "D34" '1.21'
Set-TextValue "E34" '  -8.77%  '
# Row 35
Set-TextValue "D35" '1.00'
Set-TextValue "E35" '  -0.05%  '
# Row 36
Set-TextValue "D36" '157.76'
Set-TextValue "E36" '  +0.81%  '
# Row 37
Set-TextValue "E37" '  -10.29%  '
# Row 38
Set-TextValue "D38" '18.55'
Set-TextValue "E38" '  -0.74%  '
# Row 39
Set-TextValue "D39" '18.40'
Set-TextValue "E39" '  -1.96%  '
# Row 40
Set-TextValue "D40" '1.33'
Set-TextValue "E40" '  -7.97%  '
# Row 41
Set-TextValue "E41" '  -0.04%  '
# Row 42
Set-TextValue "E42" '  -8.26%  '
# Row 43
Set-TextValue "D43" '0.324'
Set-TextValue "E43" '  -7.76%  '
# Row 44
Set-TextValue "D44" '4.70'
Set-TextValue "E44" '  -8.98%  '
# Row 45
Set-TextValue "D45" '38.61'
Set-TextValue "E45" '  -3.19%  '
# Row 46
Set-TextValue "D46" '2.28'
Set-TextValue "E46" '  -9.02%  '
# Row 47
Set-TextValue "D47" '140.67'
Set-TextValue "E47" '  -4.81%  '
# Row 48
Set-TextValue "D48" '3.46'
Set-TextValue "E48" '  -6.83%  '
# Row 49
Set-TextValue "D49" '0.508'
Set-TextValue "E49" '  -8.54%  '
# Row 50
Set-TextValue "D50" '1.58'
Set-TextValue "E50" '  -8.95%  '
# Row 51
Set-TextValue "D51" '0.0728'
Set-TextValue "E51" '  -3.13%  '
